$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8 through 19 (old rows beyond the new extent of the table)
$ws.Range("A8:C19").EntireRow.Delete()

# Update remaining data rows 2-7 with new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'Table 1.1 Incidence of AEs by Preferred Term and System Organ Class'
$ws.Range("C2").Value = '(1) Table 1.1 Incidence of AEs by Preferred Term and System Organ Class | (2) System Organ Preferred Term Treatment Group 1 Treatment Group 2 | (3) Class (N=X) (N=Y) | (4) Any System Organ | (5) Class | (6) Any Preferred Term n (%) n (%) | (7) Blood and | (8) Lymphatic System | (9) Disorders | (10) Anemia n (%) n (%) | (11) Neutropenia n (%) n (%) | (12) Gastrointestinal | (13) Disorders | (14) Diarrhea n (%) n (%) | (15) Nausea n (%) n (%) | (16) Vomiting n (%) n (%) | (17) General Disorders | (18) and Administration | (19) Site Conditions | (20) Fatigue n (%) n (%) | (21) Pyrexia n (%) n (%) | (22) Nervous System | (23) Disorders | (24) Headache n (%) n (%) | (25) Dizziness n (%) n (%) | (26) Skin and | (27) Subcutaneous | (28) Tissue Disorders | (29) Rash n (%) n (%)'

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'Table 1.2 Incidence of Serious AEs by Preferred Term and System Organ Class'
$ws.Range("C3").Value = '(1) AE: Adverse event | (2) Adverse events before treatment start date are not considered. | (3) Adverse events before treatment start date are not considered. | (4) AE’s also include treatment emergent events. | (5) Table 1.2 Incidence of Serious AEs by Preferred Term and System Organ Class | (6) System Organ Preferred Term Treatment Group 1 Treatment Group 2 | (7) System Organ Preferred Term Treatment Group 1 Treatment Group 2 | (8) Class (N=X) (N=Y) | (9) Any System Organ | (10) Class | (11) Class | (12) Any Preferred Term n (%) n (%) | (13) Blood and | (14) Lymphatic System | (15) Lymphatic System | (16) Disorders | (17) Anemia n (%) n (%) | (18) Neutropenia n (%) n (%) | (19) Gastrointestinal | (20) Gastrointestinal | (21) Disorders | (22) Diarrhea n (%) n (%) | (23) Nausea n (%) n (%) | (24) Vomiting n (%) n (%) | (25) General Disorders | (26) and Administration | (27) Site Conditions | (28) Site Conditions | (29) Fatigue n (%) n (%) | (30) Pyrexia n (%) n (%) | (31) Nervous System | (32) Nervous System | (33) Disorders | (34) Headache n (%) n (%) | (35) Dizziness n (%) n (%)'

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'Table 1.3 Incidence of AEs Leading to Discontinuation by Preferred Term and System'
$ws.Range("C4").Value = '(1) Subcutaneous | (2) Tissue Disorders | (3) Tissue Disorders | (4) Rash n (%) n (%) | (5) Rash n (%) n (%) | (6) Pruritus n (%) n (%) | (7) AE: Adverse event | (8) AE: Adverse event | (9) Serious AE: Serious Adverse event | (10) Adverse events before treatment start date are not considered. | (11) AE’s also include treatment emergent events. | (12) Table 1.3 Incidence of AEs Leading to Discontinuation by Preferred Term and System | (13) Table 1.3 Incidence of AEs Leading to Discontinuation by Preferred Term and System | (14) Organ Class | (15) System Organ Preferred Term Treatment Group 1 Treatment Group 2 | (16) Class (N=X) (N=Y) | (17) Any System Organ | (18) Class | (19) Any Preferred Term n (%) n (%) | (20) Blood and | (21) Blood and | (22) Lymphatic System | (23) Disorders | (24) Disorders | (25) Anemia n (%) n (%) | (26) Neutropenia n (%) n (%) | (27) Gastrointestinal | (28) Disorders | (29) Diarrhea n (%) n (%) | (30) Nausea n (%) n (%) | (31) Vomiting n (%) n (%) | (32) General Disorders | (33) General Disorders | (34) and Administration | (35) Site Conditions | (36) Site Conditions | (37) Fatigue n (%) n (%) | (38) Fatigue n (%) n (%) | (39) Pyrexia n (%) n (%)'

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'Table 2.1 Summary of Laboratory Values - Hematology'
$ws.Range("C5").Value = '(1) Nervous System | (2) Disorders | (3) Headache n (%) n (%) | (4) Headache n (%) n (%) | (5) Dizziness n (%) n (%) | (6) Skin and | (7) Subcutaneous | (8) Subcutaneous | (9) Tissue Disorders | (10) Rash n (%) n (%) | (11) Pruritus n (%) n (%) | (12) Pruritus n (%) n (%) | (13) AE: Adverse event | (14) Adverse events leading to discontinuation | (15) Adverse events leading to discontinuation | (16) Adverse events before treatment start date are not considered. | (17) AE’s also include treatment emergent events. | (18) Table 2.1 Summary of Laboratory Values - Hematology | (19) Parameter Statistic Treatment Group 1 Treatment Group 2 | (20) (N=X) (N=Y) | (21) (N=X) (N=Y) | (22) Hemoglobin (g/dL) N | (23) Hemoglobin (g/dL) N | (24) Mean (SD) | (25) Median | (26) Min, Max | (27) White Blood Cell N | (28) Count (x10^9/L) | (29) Count (x10^9/L) | (30) Count (x10^9/L) | (31) Mean (SD) | (32) Mean (SD) | (33) Median | (34) Min, Max | (35) Platelet Count N | (36) (x10^9/L) | (37) (x10^9/L) | (38) Mean (SD) | (39) Mean (SD)'

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'Table 2.2 Summary of Laboratory Values - Chemistry'
$ws.Range("C6").Value = '(1) Median | (2) Min, Max | (3) Summary statistics for laboratory parameters. | (4) Values presented at baseline and each post-baseline visit. | (5) Values presented at baseline and each post-baseline visit. | (6) Table 2.2 Summary of Laboratory Values - Chemistry | (7) Table 2.2 Summary of Laboratory Values - Chemistry | (8) Parameter Statistic Treatment Group 1 Treatment Group 2 | (9) (N=X) (N=Y) | (10) Alanine N | (11) Aminotransferase | (12) (U/L) | (13) Mean (SD) | (14) Median | (15) Min, Max | (16) Aspartate N | (17) Aminotransferase | (18) (U/L) | (19) Mean (SD) | (20) Median | (21) Min, Max | (22) Total Bilirubin N | (23) (µmol/L) | (24) Mean (SD) | (25) Median | (26) Min, Max | (27) Creatinine (µmol/L) N | (28) Mean (SD) | (29) Mean (SD) | (30) Median'

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 'Table 2.3 Summary of Laboratory Values - Urinalysis'
$ws.Range("C7").Value = '(1) (N=X) (N=Y) | (2) Protein N | (3) (Qualitative) | (4) Negative n (%) n (%) | (5) Trace | (6) Figure 2.2 Kaplan-Meier Plot for Time to Discontinuation Due to Adverse Event | (7) 1.0 +-----------------------------------------------------+ | (8) | | | (9) | | | (10) | +---------Treatment Group 1 | | (11) | / | | (12) | / | | (13) | / | | (14) | / | | (15) | +-----------Treatment Group 2 | | (16) | / | | (17) | / | | (18) | / | | (19) 0.0 +-----------------------------------------------------+ | (20) 0.0 +-----------------------------------------------------+ | (21) 0 Time (Days) Max Time | (22) X-axis: Time (Days) | (23) X-axis: Time (Days) | (24) Y-axis: Probability of No Discontinuation | (25) Description: Kaplan-Meier plot showing the probability of not | (26) discontinuing due to an adverse event over time for each treatment | (27) group. | (28) group. | (29) Figure 2.3 Kaplan-Meier Plot for Overall Survival | (30) 1.0 +-----------------------------------------------------+ | (31) 1.0 +-----------------------------------------------------+ | (32) | | | (33) | | | (34) | +---------Treatment Group 1 | | (35) | / | | (36) | / | | (37) | / | | (38) | / | | (39) | / |'
